$wb = $excel.ActiveWorkbook

# --- "data" sheet: add local_particles / world_particles columns (D, E) ---
$data = $wb.Worksheets.Item("data")

$data.Range("D1").Value = "local_particles"
$data.Range("D2").Value = "sphere10_DT_local.particles"
$data.Range("E1").Value = "world_particles"
$data.Range("E2").Value = "sphere10_DT_world.particles"

$data.Range("D3").Value = "sphere20_DT_local.particles"
$data.Range("E3").Value = "sphere20_DT_world.particles"

$data.Range("D4").Value = "sphere30_DT_local.particles"
$data.Range("E4").Value = "sphere30_DT_world.particles"

$data.Range("D5").Value = "sphere40_DT_local.particles"
$data.Range("E5").Value = "sphere40_DT_world.particles"

# --- "optimize" sheet: update initial_relative_weighting value and selection ---
$optimize = $wb.Worksheets.Item("optimize")

$optimize.Range("B3").Value = 0.05

$optimize.Activate()
[void]$optimize.Range("B3").Select()
